# Apply case-normalization (UPPERCASE -> lowercase) edits to the OKVED
# description column (A2:A10), fix a small typo in the "лесное хозяйство"
# row, change the sheet zoom and active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "сельское хозяйство, охота и предоставление услуг в этих  областях"
$ws.Range("A3").Value = "лесное хозяйство, лесозаготовки и предоставление услуг в этой области"
$ws.Range("A4").Value = "рыболовство, рыбоводство и предоставление услуг в этих  областях"
$ws.Range("A5").Value = "добыча каменного угля, бурого угля и торфа"
$ws.Range("A6").Value = "добыча сырой нефти и природного газа, предоставление услуг в  этих областях"
$ws.Range("A7").Value = "добыча урановой и ториевой руд"
$ws.Range("A8").Value = "добыча металлических руд"
$ws.Range("A9").Value = "добыча прочих полезных ископаемых"
$ws.Range("A10").Value = "производство пищевых продуктов, включая напитки"

# Let Excel re-flow (autofit) the wrapped-text row heights now that the
# text has changed length/case.
$ws.Rows("2:10").EntireRow.AutoFit()

# Update view: zoom and active cell/selection.
$excel.ActiveWindow.Zoom = 130
$ws.Range("A13").Select() | Out-Null
